# Update odds for the match in row 3 (Castellon - Racing Club Ferrol)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 2.3
$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 2.03
$ws.Range("W3").Value = 7.5
$ws.Range("AD3").Value = 7.5
$ws.Range("AR3").Value = 41

# Update odds for the match in row 4 (Levante - Malaga)
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.85

# Add a new match row (row 5): Boston River - Racing Montevideo
$ws.Range("A5").Value = "8A59j2KD"
$ws.Range("B5").Value = "27/11/2024"
$ws.Range("C5").Value = "16:30"
$ws.Range("D5").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E5").Value = "Boston River"
$ws.Range("F5").Value = "Racing Montevideo"
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 3.1
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.57
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 6.5
$ws.Range("X5").Value = 10
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 21
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 7.5
$ws.Range("AD5").Value = 6
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 8
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 41
$ws.Range("AN5").Value = 4.33
$ws.Range("AO5").Value = 13
$ws.Range("AP5").Value = 26
$ws.Range("AQ5").Value = 51
$ws.Range("AR5").Value = 81
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 8.5
$ws.Range("AV5").Value = 67
# AW5 is left blank (empty cell in the source data)
$ws.Range("AX5").Value = 5
$ws.Range("AY5").Value = 19
$ws.Range("AZ5").Value = 34
$ws.Range("BA5").Value = 67
$ws.Range("BB5").Value = 101
$ws.Range("BC5").Value = 301
# BD5 is left blank (empty cell in the source data)

Write-Host "Edits applied"
